# Append two more weather-observation rows (rows 3 & 4) below the existing
# row 2 for station 79049004, and refresh row 2's own readings with the
# latest observation window.
#
# Columns: A=weatherStationId  B=startDate  C=endDate
#          D=minTemperature  E=maxTemperature  F=averageTemperature  G=medianTemperature

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatted row 2 (text station id, date-formatted B/C,
# 2-decimal-formatted D:G) down into rows 3 and 4 so the new rows inherit
# the exact same cell styles/number formats instead of plain defaults.
$ws.Range("A2:G2").Copy($ws.Range("A3:G3"))
$ws.Range("A2:G2").Copy($ws.Range("A4:G4"))

# Row 2 : 2024-12-12 07:41 -> 2024-12-13 07:37
$ws.Range("B2").Value = 45638.32013888889
$ws.Range("C2").Value = 45639.31736111111
$ws.Range("D2").Value = -0.7
$ws.Range("E2").Value = 5.5
$ws.Range("F2").Value = 1.8
$ws.Range("G2").Value = 1.3

# Row 3 : 2024-12-13 07:37 -> 2024-12-14 10:30
# (column A already holds the right text value/style, copied from row 2)
$ws.Range("B3").Value = 45639.31736111111
$ws.Range("C3").Value = 45640.4375
$ws.Range("D3").Value = -0.6
$ws.Range("E3").Value = 2.8
$ws.Range("F3").Value = 1.72
$ws.Range("G3").Value = 2

# Row 4 : 2024-12-14 10:30 -> 2024-12-15 08:22
# (column A already holds the right text value/style, copied from row 2)
$ws.Range("B4").Value = 45640.4375
$ws.Range("C4").Value = 45641.34861111111
$ws.Range("D4").Value = 1.5
$ws.Range("E4").Value = 7.1
$ws.Range("F4").Value = 4.51
$ws.Range("G4").Value = 4.6

# Extend the "number stored as text" ignored-error flag (originally only
# covering A1:G2) to the full new used range, A1:G4.
$fullRange = $ws.Range("A1:G4")
try {
    $fullRange.Errors.Item(9).Ignore = $true
} catch {
}
